# Cleanup old results for Frontiers submission
# Sort the LOD_peaks table (Table3) ascending by its "lod" column (col E),
# matching Data > Sort applied on the table before submitting the results.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LOD_peaks")
$lo = $ws.ListObjects.Item("Table3")

$lodColumn = $lo.ListColumns.Item("lod").Range

$lo.Sort.SortFields.Clear()
$lo.Sort.SortFields.Add($lodColumn, 0, 1) | Out-Null
$lo.Sort.Header = 1
$lo.Sort.MatchCase = $false
$lo.Sort.Orientation = 1
$lo.Sort.Apply()
